$d = $word.ActiveDocument

$replacements = @(
    @("683×6=4098", "776×2=1552"),
    @("376×3=1128", "941×3=2823"),
    @("513×2=1026", "101×7=707"),
    @("782×7=5474", "201×5=1005"),
    @("801×3=2403", "307×2=614"),
    @("250×2=500",  "552×7=3864"),
    @("714×4=2856", "467×2=934"),
    @("833×6=4998", "614×2=1228"),
    @("872×2=1744", "170×5=850"),
    @("968×2=1936", "905×3=2715"),
    @("450×5=2250", "672×4=2688"),
    @("811×5=4055", "493×8=3944"),
    @("285×4=1140", "505×8=4040"),
    @("197×5=985",  "433×5=2165"),
    @("405×2=810",  "354×7=2478"),
    @("790×5=3950", "911×9=8199"),
    @("344×7=2408", "778×4=3112"),
    @("506×8=4048", "273×8=2184"),
    @("750×2=1500", "821×5=4105"),
    @("719×3=2157", "418×4=1672"),
    @("140×2=280",  "399×4=1596"),
    @("883×7=6181", "132×2=264"),
    @("373×4=1492", "474×9=4266"),
    @("281×7=1967", "428×6=2568"),
    @("492×8=3936", "854×4=3416")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
